# Worked on temporal resolution
#
# The "Demand" sheet (EU27.Elec demand for the 2031 intertemporal model)
# used to hold a single time step (t=0 with 0 demand, t=1 with the full
# annual demand lumped into one value). This resolves the year into a
# finer temporal grid: t=1..12, each carrying the per-step demand value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demand")

# Row 3 (t=1) previously held the whole-year lumped demand; replace it
# with the per-timestep value and extend the series through t=12.
$ws.Range("B3").Value = 315578125
for ($row = 4; $row -le 14; $row++) {
    $t = $row - 2
    $ws.Cells.Item($row, 1).Value = $t
    $ws.Cells.Item($row, 2).Value = 315578125
}

# Column B now holds 9-digit values - widen it to fit, matching Excel's
# own "best fit" sizing for this content.
$ws.Columns.Item(2).ColumnWidth = 9.166666666666666

# The editor ended their session focused on the Demand sheet (previously
# SupIm was the active/selected tab), with D13 as the active cell.
$ws.Activate()
$ws.Range("D13").Select() | Out-Null
